# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCells = @("D2","D3","D4","D5","D6","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.198.14"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.839.15"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "241.92"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "0.6630"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.07436"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").Value = "0.2948"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "23.22"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").Value = "0.07757"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "1.842.33"
$ws.Range("E12").Value = "  +2.74%  "
$ws.Range("D13").Value = "5.024"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "0.6722"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "83.33"
$ws.Range("E15").Value = "  -3.23%  "
$ws.Range("D16").Value = "6.160"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "0.000008641"
$ws.Range("E17").Value = "  +5.24%  "
$ws.Range("D18").Value = "29.200.58"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "2.091.68"
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("D20").Value = "227.93"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "7.159"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "0.9999"
$ws.Range("D25").Value = "159.92"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "0.1412"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "8.617"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "18.06"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "1.509"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "4.130"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").Value = "4.056"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").Value = "1.188"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "0.05331"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").Value = "0.7421"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").Value = "2.654"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "1.314.73"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").Value = "2.742"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("D41").Value = "6.422"
$ws.Range("E41").Value = "  +7.79%  "
$ws.Range("D42").Value = "0.9183"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "0.08290"
$ws.Range("E43").Value = "  +5.31%  "
$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "102.95"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "1.989.99"
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("D47").Value = "65.24"
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").Value = "0.5137"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").Value = "1.753"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "0.05846"
$ws.Range("E51").Value = "  -1.31%  "
